# Fixed naive component forecaster bug - Presentation state 11.02.
#
# For each data row (rows 2-20), a new error value is inserted at the
# front of the series (column B). All existing values in that row shift
# one column to the right (B->C, C->D, ... J->K), and the old value that
# was in the rightmost occupied cell of the row falls off (since the
# sheet only keeps columns B..K, i.e. 10 quarters of data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B..K, used to shift existing row values one place to the right.
$cols = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K")

# New value to insert into column B of each row (2..20).
$newValues = @{
    2  = 0.2071908139402366
    3  = -2.272829558532021
    4  = -1.400462404925464
    5  = -0.5144918801275712
    6  = -0.5616080510579985
    7  = 0.1417647591280393
    8  = -0.4790798465348092
    9  = 0.1916007792754515
    10 = 1.573432754301089
    11 = 0.9422837133007778
    12 = 0.0678490295623069
    13 = -0.5264228954459207
    14 = 0.8949500190880419
    15 = 0.2303995154407018
    16 = 0.4008418571243615
    17 = 0.2679782848922332
    18 = -0.5417707991668423
    19 = 0.0506862842519193
    20 = -0.1624199859130616
}

for ($row = 2; $row -le 20; $row++) {

    # Determine how many of the B..K cells in this row currently hold a value.
    $count = 0
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $cell = $ws.Range($cols[$i] + $row)
        if ($cell.Value2 -ne $null) {
            $count = $i + 1
        }
    }

    # Shift existing values one column to the right, starting from the
    # rightmost populated cell and working back towards column B so that
    # values are not clobbered before they are read. If the row is
    # already full (10 values, columns B..K), the value in K is dropped.
    $lastIndex = [Math]::Min($count, $cols.Length - 1) - 1
    for ($i = $lastIndex; $i -ge 0; $i--) {
        $srcCell = $ws.Range($cols[$i] + $row)
        $dstCell = $ws.Range($cols[$i + 1] + $row)
        $dstCell.Value2 = $srcCell.Value2
    }

    # Insert the new value at the front of the series (column B).
    $ws.Range("B" + $row).Value2 = $newValues[$row]
}
